# Trade #75 closed at 2026-02-17 21:13:10 - unknown UNKNOWN +0.000%
#
# Applies:
#   - Summary sheet KPI refresh (capital/P&L/trade counters)
#   - Strategy Status row for MarketMaking refreshed
#   - All Trades: trade #103 (row 104) flips OPEN -> CLOSED (early_exit), and a
#     brand-new open trade #136 is appended as row 137
#   - MarketMaking: same trade #103 (row 71) flips OPEN -> CLOSED (early_exit),
#     and trade #136 is appended as row 104

function Set-TextValue {
    # Writes a literal text value without Excel's COM layer re-interpreting
    # date-shaped strings ("2026-02-17") as date serials.
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1401.15             # Current Capital
$summary.Range("B4").Value = 0.9399999999999999  # Total P&L $
$summary.Range("B6").Value = 103                 # Total Trades
$summary.Range("B7").Value = 49                  # Winning Trades
$summary.Range("B9").Value = 47.57               # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status - MarketMaking row
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 101.15   # Capital
$status.Range("D5").Value = 70       # Trades
$status.Range("E5").Value = 0.83     # P&L $
$status.Range("F5").Value = 1.15     # P&L %
$status.Range("G5").Value = 50       # Win Rate %

# ---------------------------------------------------------------------------
# All Trades - trade #103 (sheet row 104) closes out via early_exit
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(104, 7).Value = 0.96           # Exit Price
$allTrades.Cells.Item(104, 8).Value = "CLOSED"       # Status
$allTrades.Cells.Item(104, 9).Value = 1.0526         # P&L %
$allTrades.Cells.Item(104, 10).Value = 0.01          # P&L $
$allTrades.Cells.Item(104, 11).Value = 101.15        # Capital After
$allTrades.Cells.Item(104, 12).Value = "early_exit"  # Exit Reason
$allTrades.Cells.Item(104, 13).Value = 0.15          # Duration (min)

# All Trades - new row for trade #136 (sheet row 137)
$allTrades.Cells.Item(137, 1).Value = 136
Set-TextValue $allTrades.Cells.Item(137, 2) "2026-02-17"
$allTrades.Cells.Item(137, 3).Value = "21:13:04"
$allTrades.Cells.Item(137, 4).Value = "MarketMaking"
$allTrades.Cells.Item(137, 5).Value = "UP"
$allTrades.Cells.Item(137, 6).Value = 0.95
# Exit Price (G137) stays blank - trade is still OPEN
$allTrades.Cells.Item(137, 8).Value = "OPEN"
$allTrades.Cells.Item(137, 9).Value = 0
$allTrades.Cells.Item(137, 10).Value = 0
$allTrades.Cells.Item(137, 11).Value = 101.1396151053151
# Exit Reason (L137) stays blank - trade is still OPEN
$allTrades.Cells.Item(137, 13).Value = 0
$allTrades.Cells.Item(137, 14).Value = 0
$allTrades.Cells.Item(137, 15).Value = 0
$allTrades.Cells.Item(137, 16).Value = 0.6
$allTrades.Cells.Item(137, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# MarketMaking strategy sheet - trade #103 (sheet row 71) closes out
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Cells.Item(71, 7).Value = 0.96            # Exit Price
$mm.Cells.Item(71, 8).Value = "CLOSED"        # Status
$mm.Cells.Item(71, 9).Value = 1.0526          # P&L %
$mm.Cells.Item(71, 10).Value = 0.01           # P&L $
$mm.Cells.Item(71, 11).Value = 101.15         # Capital After
$mm.Cells.Item(71, 16).Value = "early_exit"   # Exit Reason
$mm.Cells.Item(71, 17).Value = 0.15           # Duration (min)

# MarketMaking - new row for trade #136 (sheet row 104)
$mm.Cells.Item(104, 1).Value = 136
Set-TextValue $mm.Cells.Item(104, 2) "2026-02-17"
$mm.Cells.Item(104, 3).Value = "21:13:04"
$mm.Cells.Item(104, 4).Value = "MarketMaking"
$mm.Cells.Item(104, 5).Value = "UP"
$mm.Cells.Item(104, 6).Value = 0.95
# Exit Price (G104) stays blank - trade is still OPEN
$mm.Cells.Item(104, 8).Value = "OPEN"
$mm.Cells.Item(104, 9).Value = 0
$mm.Cells.Item(104, 10).Value = 0
$mm.Cells.Item(104, 11).Value = 101.1396151053151
$mm.Cells.Item(104, 12).Value = 0
$mm.Cells.Item(104, 13).Value = 0
$mm.Cells.Item(104, 14).Value = 0.6
$mm.Cells.Item(104, 15).Value = "Normal spread capture: 19600 bps"
# Exit Reason (P104) stays blank - trade is still OPEN
$mm.Cells.Item(104, 17).Value = 0
